$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "...au marché d'Elyria. Des habitants..." -> "...au marché de sa ville. Des habitants..."
$found1 = $d.Content.Find.Execute(
    "marché d’Elyria.", $true, $false, $false, $false, $false,
    $true, 1, $false, "marché de sa ville.", 2)
if (-not $found1) { throw "Change 1: target text not found" }

# --- Change 2 -----------------------------------------------------------
# "...en se relevant. Elle regarda..." -> "...en se retournant. Elle regarda..."
$found2 = $d.Content.Find.Execute(
    "se relevant.", $true, $false, $false, $false, $false,
    $true, 1, $false, "se retournant.", 2)
if (-not $found2) { throw "Change 2: target text not found" }

# --- Change 3 -------------------------------------------------------------
# Relocate the hidden "_GoBack" bookmark: it used to sit just before the
# final "« Anna… »" paragraph; it now belongs right after the word
# "retournant" we just inserted above. Re-adding a bookmark with the same
# name moves it (Word keeps bookmark names unique), which also removes it
# from its old location automatically.
$full = $d.Content.Text
$idx = $full.IndexOf("retournant")
if ($idx -lt 0) { throw "Change 3: 'retournant' not found after replace" }
$pos = $idx + ("retournant").Length
$rng = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $rng)
